$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("6_Tenants")
$ws.Activate()

$ws.Range("AY1").Value = "Logo Path"
$ws.Range("AY2").Value = "https://xform-stage.janeirodigital.com/assets/images/logo.png"
$ws.Range("AZ1").Value = "BK Color"
$ws.Range("AZ2").Value = "#04773f"

$ws.Range("AY1").Font.Bold = $true
$ws.Range("AZ1").Font.Bold = $true

$excel.ActiveWindow.ScrollColumn = 42
$ws.Range("AY7").Select()
